$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text (as in the source data) instead of
# being auto-converted to numbers by Excel when assigning numeric-looking strings.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "70.184.55"
$ws.Range("E2").Value = "  +1.17%  "

# Row 3
$ws.Range("D3").Value = "3.509.31"
$ws.Range("E3").Value = "  +0.23%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "603.21"
$ws.Range("E5").Value = "  -0.12%  "

# Row 6
$ws.Range("D6").Value = "176.05"
$ws.Range("E6").Value = "  +3.97%  "

# Row 7
$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  -0.99%  "

# Row 8
$ws.Range("D8").Value = "3.503.71"
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").Value = "0.193"
$ws.Range("E10").Value = "  +0.18%  "

# Row 11
$ws.Range("E11").Value = "  +8.96%  "

# Row 12
$ws.Range("D12").Value = "0.581"
$ws.Range("E12").Value = "  +0.59%  "

# Row 13
$ws.Range("D13").Value = "46.25"
$ws.Range("E13").Value = "  -1.55%  "

# Row 14
$ws.Range("D14").Value = "0.0000275"
$ws.Range("E14").Value = "  -0.61%  "

# Row 15
$ws.Range("D15").Value = "4.071.59"
$ws.Range("E15").Value = "  +0.41%  "

# Row 16
$ws.Range("D16").Value = "8.29"
$ws.Range("E16").Value = "  -0.09%  "

# Row 17
$ws.Range("D17").Value = "610.80"
$ws.Range("E17").Value = "  +0.10%  "

# Row 18
$ws.Range("D18").Value = "3.505.62"
$ws.Range("E18").Value = "  +0.33%  "

# Row 19
$ws.Range("D19").Value = "70.201.85"
$ws.Range("E19").Value = "  +1.23%  "

# Row 20
$ws.Range("E20").Value = "  +0.99%  "

# Row 21
$ws.Range("D21").Value = "17.34"
$ws.Range("E21").Value = "  +0.97%  "

# Row 22
$ws.Range("D22").Value = "0.877"
$ws.Range("E22").Value = "  +0.11%  "

# Row 23
$ws.Range("D23").Value = "8.99"
$ws.Range("E23").Value = "  -14.34%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "97.62"
$ws.Range("E24").Value = "  +2.05%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "15.55"
$ws.Range("E25").Value = "  -0.66%  "

# Row 26
$ws.Range("D26").Value = "3.71"
$ws.Range("E26").Value = "  -3.20%  "

# Row 27
$ws.Range("E27").Value = "  +0.04%  "

# Row 28
$ws.Range("D28").Value = "2.55"
$ws.Range("E28").Value = "  -1.63%  "

# Row 29
$ws.Range("D29").Value = "33.77"
$ws.Range("E29").Value = "  +2.38%  "

# Row 30
$ws.Range("D30").Value = "8.98"
$ws.Range("E30").Value = "  -2.76%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "8.03"
$ws.Range("E31").Value = "  -4.45%  "

# Row 32
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "2.96"
$ws.Range("E32").Value = "  -3.67%  "

# Row 33
$ws.Range("D33").Value = "638.55"
$ws.Range("E33").Value = "  +14.89%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "6.85"
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value = "1.27"
$ws.Range("E35").Value = "  -3.71%  "

# Row 36
$ws.Range("D36").Value = "3.56"
$ws.Range("E36").Value = "  +3.23%  "

# Row 37
$ws.Range("D37").Value = "0.0992"
$ws.Range("E37").Value = "  -1.59%  "

# Row 38
$ws.Range("D38").Value = "10.72"
$ws.Range("E38").Value = "  +0.05%  "

# Row 39
$ws.Range("D39").Value = "0.0470"
$ws.Range("E39").Value = "  +4.41%  "

# Row 40
$ws.Range("D40").Value = "56.67"

# Row 41
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("D42").Value = "0.142"
$ws.Range("E42").Value = "  +2.07%  "

# Row 43
$ws.Range("D43").Value = "3.361.50"
$ws.Range("E43").Value = "  +0.18%  "

# Row 44
$ws.Range("D44").Value = "0.0₃0738"
$ws.Range("E44").Value = "  +6.21%  "

# Row 45
$ws.Range("D45").Value = "0.307"
$ws.Range("E45").Value = "  -5.24%  "

# Row 46
$ws.Range("D46").Value = "32.16"
$ws.Range("E46").Value = "  -2.27%  "

# Row 47
$ws.Range("E47").Value = "  +1.63%  "

# Row 48
$ws.Range("D48").Value = "2.54"
$ws.Range("E48").Value = "  -3.10%  "

# Row 49
$ws.Range("E49").Value = "  +0.43%  "

# Row 50
$ws.Range("D50").Value = "133.79"
$ws.Range("E50").Value = "  -0.38%  "

# Row 51
$ws.Range("E51").Value = "  +0.00%  "
